$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 130.955829
$ws.Range("H2").Value = 392.867487
$ws.Range("I2").Value = 0.5336535908353144
$ws.Range("J2").Value = 0.5336535908353144
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 14734.28452296197
$ws.Range("R2").Value = 132608.5607066577
$ws.Range("S2").Value = 0.1747882531482239
$ws.Range("T2").Value = 0.1747882531482239

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 130.955829
$ws.Range("H3").Value = 392.867487
$ws.Range("I3").Value = 0.5336535908353144
$ws.Range("J3").Value = 0.5336535908353144
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 13922.49902972231
$ws.Range("R3").Value = 125302.4912675008
$ws.Range("S3").Value = 0.165158293303665
$ws.Range("T3").Value = 0.165158293303665

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 130.955829
$ws.Range("H4").Value = 392.867487
$ws.Range("I4").Value = 0.5336535908353144
$ws.Range("J4").Value = 0.5336535908353144
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 16329.09909355894
$ws.Range("R4").Value = 146961.8918420305
$ws.Range("S4").Value = 0.1937070443834255
$ws.Range("T4").Value = 0.1937070443834255

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 66.39541
$ws.Range("H5").Value = 199.18623
$ws.Range("I5").Value = 0.2705656497465488
$ws.Range("J5").Value = 0.2705656497465488
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 7470.37279233072
$ws.Range("R5").Value = 67233.35513097647
$ws.Range("S5").Value = 0.08861871838450289
$ws.Range("T5").Value = 0.08861871838450289

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 66.39541
$ws.Range("H6").Value = 199.18623
$ws.Range("I6").Value = 0.2705656497465488
$ws.Range("J6").Value = 0.2705656497465488
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 7058.792559001059
$ws.Range("R6").Value = 63529.13303100954
$ws.Range("S6").Value = 0.0837362695691621
$ws.Range("T6").Value = 0.08373626956916211

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 66.39541
$ws.Range("H7").Value = 199.18623
$ws.Range("I7").Value = 0.2705656497465488
$ws.Range("J7").Value = 0.2705656497465488
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 8278.953579435354
$ws.Range("R7").Value = 74510.58221491818
$ws.Range("S7").Value = 0.09821066179288385
$ws.Range("T7").Value = 0.09821066179288386

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 48.043585
$ws.Range("H8").Value = 144.130755
$ws.Range("I8").Value = 0.1957807594181367
$ws.Range("J8").Value = 0.1957807594181367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 5405.54671219032
$ws.Range("R8").Value = 48649.92040971287
$ws.Range("S8").Value = 0.06412432620412958
$ws.Range("T8").Value = 0.06412432620412957

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 48.043585
$ws.Range("H9").Value = 144.130755
$ws.Range("I9").Value = 0.1957807594181367
$ws.Range("J9").Value = 0.1957807594181367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 5107.72808400061
$ws.Range("R9").Value = 45969.55275600549
$ws.Range("S9").Value = 0.06059139607133916
$ws.Range("T9").Value = 0.06059139607133916

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 48.043585
$ws.Range("H10").Value = 144.130755
$ws.Range("I10").Value = 0.1957807594181367
$ws.Range("J10").Value = 0.1957807594181367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 5990.634141797704
$ws.Range("R10").Value = 53915.70727617933
$ws.Range("S10").Value = 0.07106503714266796
$ws.Range("T10").Value = 0.07106503714266797
